# Refresh the cryptos price/volume table (columns D and E, rows 2-51)
# to the latest scraped values, matching the GitHub Actions data refresh.
# Numeric-looking "Price" values are written with a leading apostrophe so
# Excel stores them as literal text (quote-prefixed) instead of silently
# parsing them into floating point numbers - these price cells are plain
# text in the workbook (e.g. thousand-separated "62.367.73" style), and a
# couple of rows (D5, D6, D14, ...) look like ordinary decimals which
# Excel would otherwise auto-convert to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.367.73"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "2.454.87"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'576.03"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "'143.94"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "2.453.14"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("E11").Value = "  +2.36%  "
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("D14").Value = "'26.37"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "2.900.49"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "62.244.93"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "2.452.09"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("E19").Value = "  -3.00%  "
$ws.Range("D20").Value = "'7.15"
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").Value = "'328.49"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("E23").Value = "  -6.78%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").Value = "'9.24"
$ws.Range("E26").Value = "  +1.63%  "
$ws.Range("D27").Value = "'593.46"
$ws.Range("D28").Value = "2.575.17"
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("D29").Value = "0.0₃0960"
$ws.Range("E29").Value = "  -3.68%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").Value = "'1.43"
$ws.Range("E31").Value = "  -4.17%  "
$ws.Range("D32").Value = "'8.02"
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("D33").Value = "'1.90"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("D35").Value = "'4.93"
$ws.Range("E35").Value = "  -3.98%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("E37").Value = "  -3.61%  "
$ws.Range("D38").Value = "'0.378"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").Value = "'151.88"
$ws.Range("E39").Value = "  +3.63%  "
$ws.Range("D40").Value = "'5.37"
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("D41").Value = "'18.43"
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("E42").Value = "  -1.99%  "
$ws.Range("D43").Value = "'42.67"
$ws.Range("E43").Value = "  +1.21%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("D46").Value = "'142.54"
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("E47").Value = "  -2.83%  "
$ws.Range("D48").Value = "'0.605"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("D49").Value = "'0.0523"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").Value = "0.0₆0247"
$ws.Range("E50").Value = "  +13.80%  "
$ws.Range("D51").Value = "'19.82"
$ws.Range("E51").Value = "  -4.34%  "


# Drop the quote-prefix styling picked up from the apostrophe above so the
# cells end up with their original (default) style, matching the source
# workbook's plain un-styled text cells.
foreach ($addr in @("D5","D6","D14","D20","D21","D26","D27","D31","D32","D33","D35","D38","D39","D40","D41","D43","D46","D48","D49","D51")) {
    $ws.Range($addr).Style = "Normal"
}
